# "Add files via upload" — reshuffle the sprint backlog tabs:
#   - insert a brand-new "backlog da sprint 3º" sheet as the 2nd tab
#     (built from the 4 rows that were selected/copied out of
#     "backlog do produto", rows 10:13)
#   - update the task text in "backlog do produto" row 10 (col C) that
#     those copied rows reference, from the old placeholder text to the
#     real task description
#   - make the new sprint-3 sheet the active tab

$wb = $excel.ActiveWorkbook
$produto = $wb.Worksheets.Item("backlog do produto")

# 1) Update the task description in "backlog do produto" (row 10, col C).
#    This is the edit that later gets carried over (copied) into the new
#    sprint-3 sheet below.
$produto.Range("C10").Value = "Fazer executavel da automatização da clonagem de todos os gits"

# 2) Insert the new sprint-3 sheet right after "backlog do produto" so it
#    becomes the 2nd tab.
$sprint3 = $wb.Worksheets.Add($null, $produto)
$sprint3.Name = "backlog da sprint 3º"

# 3) Populate it with the 4 rows that were copied out of
#    "backlog do produto" (rows 10:13), renumbering column B to 1..4.
$rows = @(
    @("Desenvolvimento ", 1, "Fazer executavel da automatização da clonagem de todos os gits", "Entrega  do processo de clonagem autonoma", "média", 3),
    @("Desenvolvimento ", 2, "Acesso do cliente ", "Acesso do cliente ", "média", 3),
    @("Desenvolvimento ", 3, "Aba de pesquisa", "Criação de abas de pesquisas", "média", 3),
    @("Desenvolvimento ", 4, "Atualização do site", "Atualização do site", "média", 3)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $data = $rows[$i]
    $sprint3.Cells.Item($r, 1).Value = $data[0]
    $sprint3.Cells.Item($r, 2).Value = $data[1]
    $sprint3.Cells.Item($r, 3).Value = $data[2]
    $sprint3.Cells.Item($r, 4).Value = $data[3]
    $sprint3.Cells.Item($r, 5).Value = $data[4]
    $sprint3.Cells.Item($r, 6).Value = $data[5]
}

# 4) The new sprint-3 tab is the one left active/selected.
$sprint3.Activate()
$sprint3.Range("B4").Select()
